$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep text (string) storage to match original inline-string cells,
# since several values look numeric (e.g. "1.00", "8.80") and Excel would
# otherwise silently coerce them to numbers and drop formatting like trailing zeros.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '89.800.30'
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.079.59'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.47%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.18%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.97'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +9.22%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '618.27'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.50%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.05'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -10.46%  '
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.96%  '
# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.12%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.079.06'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.46%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.715'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.72%  '
# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.49%  '
# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.39%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.36'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.79%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.741.58'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.11%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.38'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -6.45%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.669.82'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.94%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.077.08'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.53%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.83'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.27%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000213'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.39%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.83'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.73%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '433.97'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -8.84%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.43'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +5.00%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.80'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.24%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.75'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.04%  '
# Row 26
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '86.57'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -9.94%  '
# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.78'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.89%  '
# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.67%  '
# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.02%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.10'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.58%  '
# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.07%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.156'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.81%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.194'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -7.31%  '
# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.153'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.29%  '
# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.64'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -7.57%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.73'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.05%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.13'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.50%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '496.33'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.62%  '
# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.20%  '
# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.67%  '
# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0886'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.25%  '
# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'MantraDAO'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.60'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +54.80%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.09'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.59%  '
# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.03%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.397'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.43%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '151.70'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.49%  '
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.25%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.676'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -8.69%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.36'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.74%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.31'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.07%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.12%  '
